$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Select the range C2:C26 and set the Runmode column value to "N" for rows 2-15
# (rows 16-26 already have "N"); this effectively restricts execution to the
# single remaining "Y" row (TestCase_B26 on row 27), i.e. running B suite only.
$range = $ws.Range("C2:C26")
$range.Select()

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = "N"
}
